# Commit: "Update PHPExcel version (commit 6)"
#
# The worksheet is renamed (which also re-points the hidden
# `_xlnm._FilterDatabase` defined name at the new sheet name), the header
# row is restyled to the document's default/normal font size (11pt, still
# bold), and the active selection moves from A1 to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet "Worksheet" -> "Sim, ple".
# (Excel automatically repoints the sheet-scoped _FilterDatabase defined
# name, and the autoFilter's sheet-qualified reference, at the new name.)
$ws.Name = "Sim, ple"

# Header row (A1:C1): keep bold, move off the legacy 10pt size onto the
# workbook's normal 11pt size.
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Size = 11
$headerRange.Font.Bold = $true

# Active cell moves from A1 to C1.
$ws.Range("C1").Select()
